# UndoRedoActivityDiagram.pptx edit
#
# 1. Refresh the cached "datetimeFigureOut" field text on every slide
#    layout's Date placeholder (6/7/2018 -> 11/11/2018).
# 2. "[command commits address book]" -> "[command commits model]" and
#    center the paragraph, in the activity-diagram slide.
# 3. Collapse the 3-run "Purge redundant states and then save address
#    book to addressBookStateList " label into a single run reading
#    "Purge redundant states and then save states in …StateList".

$p = $ppt.ActivePresentation

# --- 1. Date placeholders on every slide layout -----------------------
$master = $p.SlideMaster
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shape = $layout.Shapes.Item($si)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "6/7/2018") {
                $tr.Text = "11/11/2018"
            }
        }
    }
}

# --- 2 & 3. Activity-diagram slide text tweaks -------------------------
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if (-not $shape.HasTextFrame) { continue }
    $tr = $shape.TextFrame.TextRange

    if ($shape.Name -eq "TextBox 47") {
        # "[command commits address book]" -> "[command commits model]"
        $len = $tr.Length
        $rest = $tr.Characters(2, $len - 1)
        $rest.Text = "command commits model]"
        $tr.ParagraphFormat.Alignment = 2   # ppAlignCenter
        # The shape is spAutoFit; shortening the text nudges the
        # engine's recomputed height off the original cached value by
        # a font-metrics rounding hair. Put it back so only the text
        # (and alignment) actually change.
        $shape.Height = 50.9124
    }
    elseif ($shape.Name -eq "Rounded Rectangle 50") {
        # Merge the three runs into a single run with new wording.
        $len = $tr.Length
        $whole = $tr.Characters(1, $len)
        $whole.Text = "Purge redundant states and then save states in " + [char]0x2026 + "StateList"
    }
}
